$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 242, pushing the existing rows
# (242..347) down to (244..349). Inserting via a 2-row range in one shot
# keeps the date number-format (style index 2) that Excel copies from the
# row immediately below on each new row.
$ws.Range("A242:A243").EntireRow.Insert()

# New row 242: Ajo Chino Primera entry dated 2022-10-05 (serial 44839)
$ws.Cells.Item(242, 1).Value = 8
$ws.Cells.Item(242, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(242, 3).Value = "Coquimbo"
$ws.Cells.Item(242, 4).Value = 44839
$ws.Cells.Item(242, 5).Value = 4
$ws.Cells.Item(242, 6).Value = 100112003
$ws.Cells.Item(242, 7).Value = "Ajo"
$ws.Cells.Item(242, 8).Value = "Chino"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 500
$ws.Cells.Item(242, 11).Value = 19000
$ws.Cells.Item(242, 12).Value = 20000
$ws.Cells.Item(242, 13).Value = 19500
$ws.Cells.Item(242, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(242, 15).Value = "China"
$ws.Cells.Item(242, 16).Value = 1950
$ws.Cells.Item(242, 17).Value = 10
$ws.Cells.Item(242, 18).Value = "Hortaliza"

# New row 243: Ajo Chino Primera entry also dated 2022-10-05 (serial 44839)
$ws.Cells.Item(243, 1).Value = 8
$ws.Cells.Item(243, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(243, 3).Value = "Coquimbo"
$ws.Cells.Item(243, 4).Value = 44839
$ws.Cells.Item(243, 5).Value = 4
$ws.Cells.Item(243, 6).Value = 100112003
$ws.Cells.Item(243, 7).Value = "Ajo"
$ws.Cells.Item(243, 8).Value = "Chino"
$ws.Cells.Item(243, 9).Value = "Primera"
$ws.Cells.Item(243, 10).Value = 400
$ws.Cells.Item(243, 11).Value = 21000
$ws.Cells.Item(243, 12).Value = 22000
$ws.Cells.Item(243, 13).Value = 21500
$ws.Cells.Item(243, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(243, 15).Value = "China"
$ws.Cells.Item(243, 16).Value = 2150
$ws.Cells.Item(243, 17).Value = 10
$ws.Cells.Item(243, 18).Value = "Hortaliza"
